$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 18.31647966666667
$ws.Range("H2").Value = 54.949439
$ws.Range("I2").Value = 0.005487334033884006
$ws.Range("J2").Value = 0.005487334033884005
$ws.Range("M2").Value = 0.152959
$ws.Range("N2").Value = 0.458877
$ws.Range("O2").Value = 0.004761500378002596
$ws.Range("P2").Value = 0.004761500378002596
$ws.Range("Q2").Value = 2.801670413333667
$ws.Range("R2").Value = 25.215033720003
$ws.Range("S2").Value = 0.0000261279430765652
$ws.Range("T2").Value = 0.0000261279430765652
$ws.Range("G3").Value = 18.31647966666667
$ws.Range("H3").Value = 54.949439
$ws.Range("I3").Value = 0.005487334033884006
$ws.Range("J3").Value = 0.005487334033884005
$ws.Range("O3").Value = 0.9837878817404418
$ws.Range("P3").Value = 0.9837878817404418
$ws.Range("Q3").Value = 578.8615315462007
$ws.Range("R3").Value = 5209.753783915806
$ws.Range("S3").Value = 0.00539837272559698
$ws.Range("T3").Value = 0.005398372725596979
$ws.Range("G4").Value = 18.31647966666667
$ws.Range("H4").Value = 54.949439
$ws.Range("I4").Value = 0.005487334033884006
$ws.Range("J4").Value = 0.005487334033884005
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2671263333333334
$ws.Range("N4").Value = 0.8013790000000001
$ws.Range("O4").Value = 0.008315444904458803
$ws.Range("P4").Value = 0.008315444904458805
$ws.Range("Q4").Value = 4.892814052931223
$ws.Range("R4").Value = 44.035326476381
$ws.Range("S4").Value = 0.00004562962383112413
$ws.Range("T4").Value = 0.00004562962383112413
$ws.Range("G5").Value = 18.31647966666667
$ws.Range("H5").Value = 54.949439
$ws.Range("I5").Value = 0.005487334033884006
$ws.Range("J5").Value = 0.005487334033884005
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1007146666666667
$ws.Range("N5").Value = 0.302144
$ws.Range("O5").Value = 0.00313517297709673
$ws.Range("P5").Value = 0.00313517297709673
$ws.Range("Q5").Value = 1.844738144135111
$ws.Range("R5").Value = 16.602643297216
$ws.Range("S5").Value = 0.00001720374137933633
$ws.Range("T5").Value = 0.00001720374137933632
$ws.Range("I6").Value = 0.9472399998689139
$ws.Range("J6").Value = 0.9472399998689137
$ws.Range("M6").Value = 0.152959
$ws.Range("N6").Value = 0.458877
$ws.Range("O6").Value = 0.004761500378002596
$ws.Range("P6").Value = 0.004761500378002596
$ws.Range("Q6").Value = 483.6327195631811
$ws.Range("R6").Value = 4352.694476068629
$ws.Range("S6").Value = 0.004510283617435013
$ws.Range("T6").Value = 0.004510283617435012
$ws.Range("I7").Value = 0.9472399998689139
$ws.Range("J7").Value = 0.9472399998689137
$ws.Range("O7").Value = 0.9837878817404418
$ws.Range("P7").Value = 0.9837878817404418
$ws.Range("S7").Value = 0.9318832329708553
$ws.Range("T7").Value = 0.9318832329708551
$ws.Range("I8").Value = 0.9472399998689139
$ws.Range("J8").Value = 0.9472399998689137
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2671263333333334
$ws.Range("N8").Value = 0.8013790000000001
$ws.Range("O8").Value = 0.008315444904458803
$ws.Range("P8").Value = 0.008315444904458805
$ws.Range("Q8").Value = 844.6121840293205
$ws.Range("R8").Value = 7601.509656263884
$ws.Range("S8").Value = 0.007876722030209518
$ws.Range("T8").Value = 0.007876722030209518
$ws.Range("I9").Value = 0.9472399998689139
$ws.Range("J9").Value = 0.9472399998689137
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1007146666666667
$ws.Range("N9").Value = 0.302144
$ws.Range("O9").Value = 0.00313517297709673
$ws.Range("P9").Value = 0.00313517297709673
$ws.Range("Q9").Value = 318.4442114546986
$ws.Range("R9").Value = 2865.997903092288
$ws.Range("S9").Value = 0.002969761250414129
$ws.Range("T9").Value = 0.002969761250414129
$ws.Range("G10").Value = 155.6514383333333
$ws.Range("H10").Value = 466.954315
$ws.Range("I10").Value = 0.04663076369111781
$ws.Range("J10").Value = 0.0466307636911178
$ws.Range("M10").Value = 0.152959
$ws.Range("N10").Value = 0.458877
$ws.Range("O10").Value = 0.004761500378002596
$ws.Range("P10").Value = 0.004761500378002596
$ws.Range("Q10").Value = 23.80828835602834
$ws.Range("R10").Value = 214.274595204255
$ws.Range("S10").Value = 0.0002220323989418072
$ws.Range("T10").Value = 0.0002220323989418072
$ws.Range("G11").Value = 155.6514383333333
$ws.Range("H11").Value = 466.954315
$ws.Range("I11").Value = 0.04663076369111781
$ws.Range("J11").Value = 0.0466307636911178
$ws.Range("O11").Value = 0.9837878817404418
$ws.Range("P11").Value = 0.9837878817404418
$ws.Range("Q11").Value = 4919.101902805724
$ws.Range("R11").Value = 44271.91712525152
$ws.Range("S11").Value = 0.0458747802356239
$ws.Range("T11").Value = 0.04587478023562389
$ws.Range("G12").Value = 155.6514383333333
$ws.Range("H12").Value = 466.954315
$ws.Range("I12").Value = 0.04663076369111781
$ws.Range("J12").Value = 0.0466307636911178
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.2671263333333334
$ws.Range("N12").Value = 0.8013790000000001
$ws.Range("O12").Value = 0.008315444904458803
$ws.Range("P12").Value = 0.008315444904458805
$ws.Range("Q12").Value = 41.57859800004278
$ws.Range("R12").Value = 374.2073820003851
$ws.Range("S12").Value = 0.0003877555463263282
$ws.Range("T12").Value = 0.0003877555463263281
$ws.Range("G13").Value = 155.6514383333333
$ws.Range("H13").Value = 466.954315
$ws.Range("I13").Value = 0.04663076369111781
$ws.Range("J13").Value = 0.0466307636911178
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1007146666666667
$ws.Range("N13").Value = 0.302144
$ws.Range("O13").Value = 0.00313517297709673
$ws.Range("P13").Value = 0.00313517297709673
$ws.Range("Q13").Value = 15.67638272792889
$ws.Range("R13").Value = 141.08744455136
$ws.Range("S13").Value = 0.0001461955102257759
$ws.Range("T13").Value = 0.0001461955102257759
$ws.Range("G14").Value = 2.142642
$ws.Range("H14").Value = 6.427926
$ws.Range("I14").Value = 0.0006419024060843985
$ws.Range("J14").Value = 0.0006419024060843984
$ws.Range("M14").Value = 0.152959
$ws.Range("N14").Value = 0.458877
$ws.Range("O14").Value = 0.004761500378002596
$ws.Range("P14").Value = 0.004761500378002596
$ws.Range("Q14").Value = 0.327736377678
$ws.Range("R14").Value = 2.949627399102
$ws.Range("S14").Value = 0.00000305641854921164
$ws.Range("T14").Value = 0.000003056418549211639
$ws.Range("G15").Value = 2.142642
$ws.Range("H15").Value = 6.427926
$ws.Range("I15").Value = 0.0006419024060843985
$ws.Range("J15").Value = 0.0006419024060843984
$ws.Range("O15").Value = 0.9837878817404418
$ws.Range("P15").Value = 0.9837878817404418
$ws.Range("Q15").Value = 67.714596486156
$ws.Range("R15").Value = 609.4313683754041
$ws.Range("S15").Value = 0.0006314958083658633
$ws.Range("T15").Value = 0.0006314958083658632
$ws.Range("G16").Value = 2.142642
$ws.Range("H16").Value = 6.427926
$ws.Range("I16").Value = 0.0006419024060843985
$ws.Range("J16").Value = 0.0006419024060843984
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2671263333333334
$ws.Range("N16").Value = 0.8013790000000001
$ws.Range("O16").Value = 0.008315444904458803
$ws.Range("P16").Value = 0.008315444904458805
$ws.Range("Q16").Value = 0.572356101106
$ws.Range("R16").Value = 5.151204909954001
$ws.Range("S16").Value = 0.000005337704091834357
$ws.Range("T16").Value = 0.000005337704091834357
$ws.Range("G17").Value = 2.142642
$ws.Range("H17").Value = 6.427926
$ws.Range("I17").Value = 0.0006419024060843985
$ws.Range("J17").Value = 0.0006419024060843984
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1007146666666667
$ws.Range("N17").Value = 0.302144
$ws.Range("O17").Value = 0.00313517297709673
$ws.Range("P17").Value = 0.00313517297709673
$ws.Range("Q17").Value = 0.215795474816
$ws.Range("R17").Value = 1.942159273344
$ws.Range("S17").Value = 0.000002012475077489178
$ws.Range("T17").Value = 0.000002012475077489178
